# Weekly update: a new price-report row for "Ajo" (Feria Lagunitas de
# Puerto Montt) is prepended above the existing historical series, which
# pushes every existing row (473..523) down by one (-> 474..524).
#
# Insert a new blank row at 473 (this shifts rows 473-523 down to 474-524,
# and Excel automatically extends the used range / dimension to R524).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(473).Insert()

# Populate the newly inserted row with the latest week's data point.
$ws.Range("A473").Value2 = 4
$ws.Range("B473").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C473").Value2 = "Los Lagos"
$ws.Range("D473").Value2 = 45194
$ws.Range("E473").Value2 = 10
$ws.Range("F473").Value2 = 100112003
$ws.Range("G473").Value2 = "Ajo"
$ws.Range("H473").Value2 = "Chino"
$ws.Range("I473").Value2 = "Primera"
$ws.Range("J473").Value2 = 80
$ws.Range("K473").Value2 = 24000
$ws.Range("L473").Value2 = 24000
$ws.Range("M473").Value2 = 24000
$ws.Range("N473").Value2 = "`$/caja 10 kilos"
$ws.Range("O473").Value2 = "China"
$ws.Range("P473").Value2 = 2400
$ws.Range("Q473").Value2 = 10
$ws.Range("R473").Value2 = "Hortaliza"

# Ensure the date column keeps the same date number format used by the
# rest of the column (inherited from the Insert, but set explicitly too).
$ws.Range("D473").NumberFormat = $ws.Range("D474").NumberFormat
